# "implemented my own wait method"
#
# The `login` sheet used a column of "PASS" markers (column G, rows 2-9)
# that mirrored a per-row sleep-based wait step. Those markers are no
# longer needed, so they are cleared out, the trailing "sleep" row (row 9)
# is repurposed to a custom "debug" step, and its stray count value (D9)
# is dropped too. The `login` sheet becomes the active tab (with G2:G12
# selected) instead of `suite`.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("login")

# Row 9 used to be a "sleep" action with a 5000ms value; replace it with
# a plain "debug" marker and drop the now-unused value/status cells.
$ws.Range("C9").Value = "debug"
$ws.Range("D9").ClearContents()

# Drop the per-row "PASS" status markers in column G (rows 2-9).
$ws.Range("G2:G9").ClearContents()

# Make "login" the active sheet/tab, with G2:G12 selected.
$ws.Activate()
$ws.Range("G2:G12").Select()
